$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (VBA/COM RGB encoding: R + G*256 + B*65536)
$red   = 255       # RGB(255,0,0)   -> fill "00FF0000" (existing style s="2")
$green = 582476     # RGB(76,227,8)  -> fill "004CE308" (existing style s="3")

function Set-TextValue($range, $value) {
    # Force the cell to keep its value as literal text, even if it looks
    # like a number or a date (mirrors typing into a "Text" formatted cell).
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ----- Row 3 (gaido martin) -----
$ws.Range("H3").Interior.Color = $red
Set-TextValue $ws.Range("H3") "2500"
Set-TextValue $ws.Range("I3") "06/02/2024"
Set-TextValue $ws.Range("J3") "06/03/2024"
$ws.Range("K3").Interior.Color = $green
$ws.Range("K3").Value = "Regular"

# ----- Row 4 (gonzales matiass) -----
$ws.Range("H4").Interior.Color = $red
Set-TextValue $ws.Range("H4") "1000"
Set-TextValue $ws.Range("I4") "06/02/2024"
Set-TextValue $ws.Range("J4") "06/03/2024"
$ws.Range("K4").Interior.Color = $green
$ws.Range("K4").Value = "Regular"

# ----- Row 8 (Rodriguez Carina) -----
$ws.Range("H8").Interior.Color = $green
Set-TextValue $ws.Range("H8") "-900"
Set-TextValue $ws.Range("I8") "06/02/2024"
Set-TextValue $ws.Range("J8") "20/02/2024"
$ws.Range("K8").Interior.Color = $green
$ws.Range("K8").Value = "Regular"
